$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add new time log entry: 09/29 (serial 45198), Internship, standard description
$ws.Range("A9").Value = 45198
$ws.Range("A9").NumberFormat = "d-mmm"
$ws.Range("B9").Value = "Internship"
$ws.Range("C9").Value = "Completed 8 hours assisting with daily operations"

$ws.Range("C10").Select()
